$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block updates
$ws.Range("C2").Value = "Hartmut"
# B3 holds a long digit string that must stay plain TEXT (not be
# auto-coerced to a number by Excel). Write it with a leading apostrophe
# so it is stored as text, then copy C3's formats on top so the cell
# keeps its original style (the apostrophe write nudges the style).
$ws.Range("B3").Value = "'2570314725427075"
$ws.Range("C3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("C3").Value = "Mohaupt"

# Opening balance
$ws.Range("D5").Value = "KONTOSTAND AM 03.04.2024"

# Transaction rows 6-8 (existing rows, values updated in place)
$ws.Range("B6").Value = "04.04."
$ws.Range("C6").Value = "05.04."
$ws.Range("D6").Value = "PAYPAL EREXGY"
$ws.Range("E6").Value = "59,67-"

$ws.Range("B7").Value = "05.04."
$ws.Range("C7").Value = "06.04."
$ws.Range("D7").Value = "MCDONALDS Luckenwalde"
$ws.Range("E7").Value = "18,44-"

$ws.Range("B8").Value = "09.04."
$ws.Range("C8").Value = "10.04."
$ws.Range("D8").Value = "PAYPAL IXHYMM"
$ws.Range("E8").Value = "37,78-"

# Rows 9-10 were previously empty - now populated with new transactions.
# Their amount cells (E9/E10) previously used the "blank tail" styles
# (s=13 / s=12); once filled in they take on the same right-aligned
# amount style as the other amount cells in the block (E6-E8), so copy
# that formatting across via PasteSpecial (values already set first).
$ws.Range("B9").Value = "13.04."
$ws.Range("C9").Value = "14.04."
$ws.Range("D9").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 84333501"
$ws.Range("E9").Value = "86,06-"
$ws.Range("E8").Copy()
$ws.Range("E9").PasteSpecial(-4122)

$ws.Range("B10").Value = "14.04."
$ws.Range("C10").Value = "15.04."
$ws.Range("D10").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E10").Value = "68,64-"
$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)

# Closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 18.04.2024"
$ws.Range("E12").Value = "270,59-"

# Next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 28.04.2024"
